$d = $word.ActiveDocument

$d.Content.Find.Execute("komportashon aki sosodé", $true, $false, $false, $false, $false,
                         $true, 1, $false, "komportashon akí sosodé", 2)

$d.Content.Find.Execute("e lo ta enfoká riba e yu ku el a skohe", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e lo ta enfokando riba e yu ku el a skohe", 2)

$d.Content.Find.Execute("e abilidatnan ku ta siña den e programa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "e abilidatnan ku e ta siña den e programa", 2)

$d.Content.Find.Execute("kon pa baha e aplikahson si nan", $true, $false, $false, $false, $false,
                         $true, 1, $false, "kon pa baha e aplikashon si nan", 2)
